$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.858.94"
$ws.Range("E2").Value = "  +6.07%  "
$ws.Range("D3").Value = "2.755.24"
$ws.Range("E3").Value = "  +4.38%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'117.37"
$ws.Range("E5").Value = "  +5.85%  "
$ws.Range("D6").Value = "'332.31"
$ws.Range("E6").Value = "  +2.85%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.575"
$ws.Range("E9").Value = "  +6.43%  "
$ws.Range("D10").Value = "'41.83"
$ws.Range("E10").Value = "  +5.43%  "
$ws.Range("D11").Value = "'0.0834"
$ws.Range("E11").Value = "  +3.02%  "
$ws.Range("D12").Value = "'20.01"
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("E13").Value = "  +2.92%  "
$ws.Range("D14").Value = "'7.61"
$ws.Range("E14").Value = "  +5.40%  "
$ws.Range("D15").Value = "3.186.45"
$ws.Range("E15").Value = "  +4.80%  "
$ws.Range("D16").Value = "2.770.79"
$ws.Range("E16").Value = "  +5.27%  "
$ws.Range("D17").Value = "'0.884"
$ws.Range("E17").Value = "  +2.22%  "
$ws.Range("D18").Value = "51.781.78"
$ws.Range("E18").Value = "  +5.97%  "
$ws.Range("D19").Value = "'13.49"
$ws.Range("E19").Value = "  +5.01%  "
$ws.Range("D20").Value = "'3.05"
$ws.Range("E20").Value = "  +5.76%  "
$ws.Range("E21").Value = "  +2.32%  "
$ws.Range("E22").Value = "  +2.31%  "
$ws.Range("D23").Value = "'278.87"
$ws.Range("E23").Value = "  +2.92%  "
$ws.Range("D24").Value = "'69.82"
$ws.Range("E24").Value = "  +0.89%  "
$ws.Range("E25").Value = "  +4.39%  "
$ws.Range("D26").Value = "'26.80"
$ws.Range("E26").Value = "  +2.70%  "
$ws.Range("D27").Value = "'4.15"
$ws.Range("E27").Value = "  +0.86%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("D29").Value = "'10.26"
$ws.Range("E29").Value = "  +1.17%  "
$ws.Range("D30").Value = "'2.23"
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("E31").Value = "  +1.72%  "
$ws.Range("D32").Value = "'35.08"
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("D33").Value = "'50.52"
$ws.Range("E33").Value = "  +2.15%  "
$ws.Range("E34").Value = "  +2.72%  "
$ws.Range("E35").Value = "  +3.10%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").Value = "'19.09"
$ws.Range("E37").Value = "  -0.98%  "
$ws.Range("E38").Value = "  +2.45%  "
$ws.Range("E39").Value = "  +1.00%  "
$ws.Range("D40").Value = "'3.18"
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("D41").Value = "'130.75"
$ws.Range("E41").Value = "  +4.49%  "
$ws.Range("D42").Value = "'0.0346"
$ws.Range("E42").Value = "  +10.58%  "
$ws.Range("D43").Value = "'23.07"
$ws.Range("E43").Value = "  +2.43%  "
$ws.Range("D44").Value = "'0.113"
$ws.Range("E44").Value = "  +2.67%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'2.46"
$ws.Range("E45").Value = "  +15.87%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'2.26"
$ws.Range("E46").Value = "  +5.39%  "
$ws.Range("D47").Value = "2.111.86"
$ws.Range("E47").Value = "  +2.09%  "
$ws.Range("D48").Value = "'3.34"
$ws.Range("E48").Value = "  +3.52%  "
$ws.Range("D49").Value = "'2.23"
$ws.Range("E49").Value = "  +2.23%  "
$ws.Range("D50").Value = "'5.58"
$ws.Range("E50").Value = "  +8.03%  "
$ws.Range("D51").Value = "'8.97"
$ws.Range("E51").Value = "  +0.19%  "
